$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same "last changed" date (serial 45180,
# i.e. 2023-09-11) for every data row (2..463). Bump it by one day to
# 45181 (2023-09-12) across the whole range in one shot.
$ws.Range("C2:C463").Value = 45181
